$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("minhkhoi") got its purchase/history data corrected ("sua loi mua do").
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 11
$ws.Range("E4").Value = 9
$ws.Range("I4").Value = ";2;2"
$ws.Range("I4").NumberFormat = $ws.Range("I7").NumberFormat
$ws.Range("F4").Value = ";0;0;0;1;1;0;0;0;0;0;0"
$ws.Range("G4").Value = ";0;0;0;0;3;42;33;31;0;42;43"
$ws.Range("H4").Value = ";-100;-200;-100;+600;+700;-200;-500;-100;-120;-100;-220"
